$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-4 with new TPM-derived values (columns G, H, M, N, O, P, Q, R, S, T)
$ws.Range("G2").Value = 0.2284785
$ws.Range("H2").Value = 0.456957
$ws.Range("M2").Value = 2.330840333333333
$ws.Range("N2").Value = 6.992521
$ws.Range("O2").Value = 0.6715345129768794
$ws.Range("P2").Value = 0.7003397275969581
$ws.Range("Q2").Value = 0.5325469030994999
$ws.Range("R2").Value = 3.195281418597
$ws.Range("S2").Value = 0.6715345129768794
$ws.Range("T2").Value = 0.7003397275969581

$ws.Range("G3").Value = 0.2284785
$ws.Range("H3").Value = 0.456957
$ws.Range("O3").Value = 0.1894790521235985
$ws.Range("P3").Value = 0.1976066831789769
$ws.Range("Q3").Value = 0.150262541181
$ws.Range("R3").Value = 0.901575247086
$ws.Range("S3").Value = 0.1894790521235985
$ws.Range("T3").Value = 0.1976066831789769

$ws.Range("G4").Value = 0.2284785
$ws.Range("H4").Value = 0.456957
$ws.Range("M4").Value = 0.42828
$ws.Range("N4").Value = 0.85656
$ws.Range("O4").Value = 0.123391035029171
$ws.Range("P4").Value = 0.08578923067523865
$ws.Range("Q4").Value = 0.09785277198
$ws.Range("R4").Value = 0.39141108792
$ws.Range("S4").Value = 0.123391035029171
$ws.Range("T4").Value = 0.08578923067523865

# Add new row 5 for Resolving-Mac target cluster
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Fgf5"
$ws.Range("C5").Value = "Fgfr3"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.5
$ws.Range("G5").Value = 0.2284785
$ws.Range("H5").Value = 0.456957
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.05413033333333334
$ws.Range("N5").Value = 0.162391
$ws.Range("O5").Value = 0.01559539987035126
$ws.Range("P5").Value = 0.01626435854882633
$ws.Range("Q5").Value = 0.0123676173645
$ws.Range("R5").Value = 0.07420570418700001
$ws.Range("S5").Value = 0.01559539987035126
$ws.Range("T5").Value = 0.01626435854882633
